$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5228
$ws1.Range("F7").Value = 440
$ws1.Range("F8").Value = 621
$ws1.Range("F9").Value = 907
$ws1.Range("F11").Value = 72
$ws1.Range("F17").Value = 1747
$ws1.Range("F18").Value = 1450
$ws1.Range("F19").Value = 818
$ws1.Range("F28").Value = 2513
$ws1.Range("F34").Value = 248
$ws1.Range("F42").Value = 45
$ws1.Range("F43").Value = 44
$ws1.Range("F44").Value = 58

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 31

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 31
$ws4.Range("F7").Value = 5228
$ws4.Range("F8").Value = 440
$ws4.Range("F9").Value = 621
$ws4.Range("F12").Value = 907
$ws4.Range("F15").Value = 72
$ws4.Range("F22").Value = 1747
$ws4.Range("F23").Value = 1450
$ws4.Range("F24").Value = 818
$ws4.Range("F33").Value = 2513
$ws4.Range("F38").Value = 248
$ws4.Range("F45").Value = 44
$ws4.Range("F46").Value = 58
